$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Extend the table with a new year column (S) for 2022, copying the
# formatting of the preceding year column (R) so the new cells match the
# existing table style.
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial($xlPasteFormats)
$ws.Range("S3").Value = 2022

$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial($xlPasteFormats)
$ws.Range("S4").Value = 0.071025550219041236

# Unify the first three column widths into one consistent width.
$ws.Range("A:C").ColumnWidth = 32.59

# Update the active selection, as recorded by the last user interaction.
$ws.Range("F14").Select() | Out-Null
